$d = $word.ActiveDocument

# The document currently ends with:
#   Para n-1: hyperlink paragraph (ilvl 1)
#   Para n  : empty paragraph holding the _GoBack bookmark (ilvl 0)
#
# We need to insert, right before the bookmark paragraph, two new list
# items:
#   1) "Library API Documentation"                           (ilvl 0)
#   2) hyperlink "http://www.digitalrune.com/Documentation/" (ilvl 1)
#      followed by a trailing space run.
#
# Note: once a Range.InsertBefore() shifts paragraphs around, previously
# fetched Paragraph objects can report stale .Index values, so every
# insertion point below is re-fetched fresh by (fixed) numeric position
# instead of relying on cached paragraph references.

$n = $d.Paragraphs.Count

# Insert a blank paragraph right before the bookmark paragraph (item n).
# It inherits the bookmark paragraph's ListParagraph / ilvl=0 / numId=2
# formatting automatically.
$d.Paragraphs.Item($n).Range.InsertBefore("`r")

# Insert the "Library API Documentation" paragraph right before the
# blank paragraph we just created (which is now item n again, since it
# got pushed down by one).
$d.Paragraphs.Item($n).Range.InsertBefore("Library API Documentation`r")

# Layout is now:
#   Para n    : "Library API Documentation"   (ilvl 0)
#   Para n+1  : "" (blank, to become the hyperlink paragraph)
#   Para n+2  : bookmark paragraph

$hlPara = $d.Paragraphs.Item($n + 1)

# Promote it to the second list level (ilvl=1 <=> ListLevelNumber=2).
$hlPara.Range.ListFormat.ListLevelNumber = 2

$url = "http://www.digitalrune.com/Documentation/"

# Type the URL text plus a trailing space, then turn just the URL part
# into a real hyperlink (leaving the trailing space as plain text),
# matching the existing hyperlink paragraph above it.
$hlPara.Range.InsertBefore($url + " ")

$hlPara = $d.Paragraphs.Item($n + 1)
$urlStart = $hlPara.Range.Start
$urlEnd = $urlStart + $url.Length
$urlRange = $d.Range($urlStart, $urlEnd)

$d.Hyperlinks.Add($urlRange, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null
